$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8387
$ws.Range("E2").Value = 2782
$ws.Range("F2").Value = 2782
$ws.Range("G2").Value = 2890
$ws.Range("H2").Value = 2275
$ws.Range("I2").Value = 2300
$ws.Range("J2").Value = -24
$ws.Range("K2").Value = 16988
$ws.Range("L2").Value = 3232
$ws.Range("M2").Value = 13756
$ws.Range("N2").Value = 13606
$ws.Range("O2").Value = 150
$ws.Range("P2").Value = 110
$ws.Range("Q2").Value = 2636
$ws.Range("R2").Value = -1755
$ws.Range("S2").Value = -123
$ws.Range("T2").Value = 195
$ws.Range("U2").Value = 2441
$ws.Range("W2").Value = 33.17
$ws.Range("X2").Value = 27.13
$ws.Range("Y2").Value = 18.38
$ws.Range("Z2").Value = 14.32
$ws.Range("AA2").Value = 23.5
$ws.Range("AB2").Value = 13284.86
$ws.Range("AC2").Value = 10487
$ws.Range("AD2").Value = 17.35
$ws.Range("AE2").Value = 68130
$ws.Range("AF2").Value = 2.67
$ws.Range("AG2").Value = 3430
$ws.Range("AH2").Value = 1.88
$ws.Range("AI2").Value = 29.79
$ws.Range("AJ2").Value = 21929022

# Row 3
$ws.Range("D3").Value = 8383
$ws.Range("E3").Value = 2375
$ws.Range("F3").Value = 2375
$ws.Range("G3").Value = 2395
$ws.Range("H3").Value = 1664
$ws.Range("I3").Value = 1654
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 22192
$ws.Range("L3").Value = 4303
$ws.Range("M3").Value = 17889
$ws.Range("N3").Value = 17695
$ws.Range("O3").Value = 194
$ws.Range("P3").Value = 110
$ws.Range("Q3").Value = 2413
$ws.Range("R3").Value = -5169
$ws.Range("S3").Value = 3102
$ws.Range("T3").Value = 177
$ws.Range("U3").Value = 2237
$ws.Range("W3").Value = 28.33
$ws.Range("X3").Value = 19.84
$ws.Range("Y3").Value = 10.57
$ws.Range("Z3").Value = 8.49
$ws.Range("AA3").Value = 24.05
$ws.Range("AB3").Value = 16156.9
$ws.Range("AC3").Value = 7542
$ws.Range("AD3").Value = 28.24
$ws.Range("AE3").Value = 81045
$ws.Range("AF3").Value = 2.63
$ws.Range("AG3").Value = 2747
$ws.Range("AH3").Value = 1.29
$ws.Range("AI3").Value = 36.32
$ws.Range("AJ3").Value = 21929022

# Row 4
$ws.Range("D4").Value = 9836
$ws.Range("E4").Value = 3288
$ws.Range("F4").Value = 3288
$ws.Range("G4").Value = 3461
$ws.Range("H4").Value = 2714
$ws.Range("I4").Value = 2723
$ws.Range("J4").Value = -9
$ws.Range("K4").Value = 23608
$ws.Range("L4").Value = 4653
$ws.Range("M4").Value = 18955
$ws.Range("N4").Value = 18856
$ws.Range("O4").Value = 99
$ws.Range("P4").Value = 110
$ws.Range("Q4").Value = 1120
$ws.Range("R4").Value = -1387
$ws.Range("S4").Value = -738
$ws.Range("T4").Value = 166
$ws.Range("U4").Value = 954
$ws.Range("V4").Value = 1497
$ws.Range("W4").Value = 33.43
$ws.Range("X4").Value = 27.59
$ws.Range("Y4").Value = 14.9
$ws.Range("Z4").Value = 11.85
$ws.Range("AA4").Value = 24.55
$ws.Range("AB4").Value = 18017.34
$ws.Range("AC4").Value = 12416
$ws.Range("AD4").Value = 19.93
$ws.Range("AE4").Value = 88772
$ws.Range("AF4").Value = 2.79
$ws.Range("AG4").Value = 3820
$ws.Range("AH4").Value = 1.54
$ws.Range("AI4").Value = 29.8
$ws.Range("AJ4").Value = 21929022

# Row 5
$ws.Range("D5").Value = 17587
$ws.Range("E5").Value = 5850
$ws.Range("F5").Value = 5850
$ws.Range("G5").Value = 6102
$ws.Range("H5").Value = 4440
$ws.Range("I5").Value = 4410
$ws.Range("J5").Value = 31
$ws.Range("K5").Value = 35266
$ws.Range("L5").Value = 7973
$ws.Range("M5").Value = 27292
$ws.Range("N5").Value = 27212
$ws.Range("O5").Value = 80
$ws.Range("P5").Value = 110
$ws.Range("Q5").Value = 5952
$ws.Range("R5").Value = -4596
$ws.Range("S5").Value = -797
$ws.Range("T5").Value = 288
$ws.Range("U5").Value = 5664
$ws.Range("V5").Value = 1498
$ws.Range("W5").Value = 33.26
$ws.Range("X5").Value = 25.25
$ws.Range("Y5").Value = 19.14
$ws.Range("Z5").Value = 15.09
$ws.Range("AA5").Value = 29.21
$ws.Range("AB5").Value = 21318.41
$ws.Range("AC5").Value = 20104
$ws.Range("AD5").Value = 22.26
$ws.Range("AE5").Value = 128056
$ws.Range("AF5").Value = 3.49
$ws.Range("AG5").Value = 7280
$ws.Range("AH5").Value = 1.63
$ws.Range("AI5").Value = 35.08
$ws.Range("AJ5").Value = 21939022

# Row 6
$ws.Range("D6").Value = 17151
$ws.Range("E6").Value = 6149
$ws.Range("F6").Value = 6149
$ws.Range("G6").Value = 6374
$ws.Range("H6").Value = 4215
$ws.Range("I6").Value = 4182
$ws.Range("K6").Value = 29413
$ws.Range("L6").Value = 5623
$ws.Range("M6").Value = 23790
$ws.Range("N6").Value = 23677
$ws.Range("P6").Value = 110
$ws.Range("Q6").Value = 3528
$ws.Range("R6").Value = 683
$ws.Range("S6").Value = -4244
$ws.Range("T6").Value = 267
$ws.Range("U6").Value = 3261
$ws.Range("V6").Value = 1551
$ws.Range("W6").Value = 35.85
$ws.Range("X6").Value = 24.57
$ws.Range("Y6").Value = 16.43
$ws.Range("Z6").Value = 13.03
$ws.Range("AA6").Value = 23.64
$ws.Range("AB6").Value = 23669.54
$ws.Range("AC6").Value = 19061
$ws.Range("AD6").Value = 24.47
$ws.Range("AE6").Value = 114980
$ws.Range("AF6").Value = 4.06
$ws.Range("AG6").Value = 6050
$ws.Range("AH6").Value = 1.3
$ws.Range("AI6").Value = 29.79
$ws.Range("AJ6").Value = 21939022

# Row 7
$ws.Range("D7").Value = 16967
$ws.Range("E7").Value = 5179
$ws.Range("G7").Value = 5919
$ws.Range("H7").Value = 4425
$ws.Range("I7").Value = 4418
$ws.Range("K7").Value = 33247
$ws.Range("L7").Value = 6487
$ws.Range("M7").Value = 26760
$ws.Range("N7").Value = 26592
$ws.Range("P7").Value = 110
$ws.Range("Q7").Value = 4814
$ws.Range("R7").Value = -758
$ws.Range("S7").Value = -1013
$ws.Range("T7").Value = 387
$ws.Range("U7").Value = 3936
$ws.Range("W7").Value = 30.52
$ws.Range("X7").Value = 26.08
$ws.Range("Y7").Value = 17.58
$ws.Range("Z7").Value = 14.12
$ws.Range("AA7").Value = 24.24
$ws.Range("AC7").Value = 20131
$ws.Range("AD7").Value = 31.59
$ws.Range("AE7").Value = 129041
$ws.Range("AF7").Value = 4.93
$ws.Range("AG7").Value = 6416
$ws.Range("AH7").Value = 1.01
$ws.Range("AI7").Value = 31.88

# Row 8
$ws.Range("D8").Value = 24881
$ws.Range("E8").Value = 10002
$ws.Range("G8").Value = 10670
$ws.Range("H8").Value = 7960
$ws.Range("I8").Value = 7997
$ws.Range("K8").Value = 40736
$ws.Range("L8").Value = 7463
$ws.Range("M8").Value = 33286
$ws.Range("N8").Value = 33110
$ws.Range("P8").Value = 110
$ws.Range("Q8").Value = 7919
$ws.Range("R8").Value = -2313
$ws.Range("S8").Value = -1720
$ws.Range("T8").Value = 319
$ws.Range("U8").Value = 8137
$ws.Range("W8").Value = 40.2
$ws.Range("X8").Value = 31.99
$ws.Range("Y8").Value = 26.79
$ws.Range("Z8").Value = 21.52
$ws.Range("AA8").Value = 22.42
$ws.Range("AC8").Value = 36427
$ws.Range("AD8").Value = 17.46
$ws.Range("AE8").Value = 160669
$ws.Range("AF8").Value = 3.96
$ws.Range("AG8").Value = 7864
$ws.Range("AH8").Value = 1.24
$ws.Range("AI8").Value = 21.59

# Row 9
$ws.Range("D9").Value = 28185
$ws.Range("E9").Value = 11674
$ws.Range("G9").Value = 12463
$ws.Range("H9").Value = 9256
$ws.Range("I9").Value = 9208
$ws.Range("K9").Value = 48297
$ws.Range("L9").Value = 7533
$ws.Range("M9").Value = 40805
$ws.Range("N9").Value = 40657
$ws.Range("P9").Value = 110
$ws.Range("Q9").Value = 9022
$ws.Range("R9").Value = -1775
$ws.Range("S9").Value = -2050
$ws.Range("T9").Value = 347
$ws.Range("U9").Value = 9537
$ws.Range("W9").Value = 41.42
$ws.Range("X9").Value = 32.84
$ws.Range("Y9").Value = 24.97
$ws.Range("Z9").Value = 20.79
$ws.Range("AA9").Value = 18.46
$ws.Range("AC9").Value = 41944
$ws.Range("AD9").Value = 15.16
$ws.Range("AE9").Value = 197292
$ws.Range("AF9").Value = 3.22
$ws.Range("AG9").Value = 8751
$ws.Range("AH9").Value = 1.38
$ws.Range("AI9").Value = 20.86

# Remove V2 and V3 entirely (cell element removed on save)
$ws.Range("V2").ClearContents()
$ws.Range("V3").ClearContents()
